$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$f2 = @"
<rpc-reply message-id="urn:uuid:6817521b-e237-416a-b55a-ef0ddf61df8a">
  <data>
    <network-instances>
      <network-instance>
        <name>Prueba_LxVPN</name>
        <config>
          <name>Prueba_LxVPN</name>
          <type>oc-ni-types:L3VRF</type>
        </config>
        <interfaces>
          <interface>
            <id>GigabitEthernet0/3/2</id>
            <config>
              <id>GigabitEthernet0/3/2</id>
              <interface>GigabitEthernet0/3/2</interface>
              <subinterface>0</subinterface>
            </config>
          </interface>
        </interfaces>
        <protocols>
          <protocol>
            <identifier>oc-pol-types:OSPF</identifier>
            <name>22</name>
            <config>
              <identifier>oc-pol-types:OSPF</identifier>
              <name>22</name>
            </config>
            <ospfv2>
              <global>
                <config>
                  <router-id>172.16.1.3</router-id>
                </config>
              </global>
            </ospfv2>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:STATIC</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:STATIC</identifier>
              <name>default</name>
            </config>
          </protocol>
          <protocol>
            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
            <name>default</name>
            <config>
              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>
              <name>default</name>
            </config>
          </protocol>
        </protocols>
      </network-instance>
    </network-instances>
  </data>
</rpc-reply>
"@

$g2 = @"
  <edit-config>
    <target>
      <candidate/>
    </target>
    <config>
      <network-instances xmlns="http://openconfig.net/yang/network-instance">
        <network-instance>
          <name>Prueba_LxVPN</name>
          <config>
            <name>Prueba_LxVPN</name>
            <type xmlns:oc-ni-types="http://openconfig.net/yang/network-instance-types">oc-ni-types:L3VRF</type>
          </config>
          <protocols>
            <protocol>
              <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:STATIC</identifier>
              <name>default</name>
              <config>
                <identifier xmlns:oc-pol-types="http://openconfig.net/yang/policy-types">oc-pol-types:STATIC</identifier>
                <name>default</name>
              </config>
              <static-routes>
                <static>
                  <prefix>192.168.100.0/24</prefix>
                  <config>
                    <prefix>192.168.100.0/24</prefix>
                  </config>
                  <next-hops>
                    <next-hop>
                      <index>10.10.10.10</index>
                      <config>
                        <index>10.10.10.10</index>
                        <metric>150</metric>
                        <next-hop>10.10.10.10</next-hop>
                      </config>
                    </next-hop>
                  </next-hops>
                </static>
              </static-routes>
            </protocol>
          </protocols>
        </network-instance>
      </network-instances>
    </config>
  </edit-config>
"@

# Remove the trailing newline that PowerShell here-strings append, the F2 cell keeps one
# trailing newline (after </rpc-reply>) while G2 has none (after </edit-config>).
$f2 = $f2 -replace "`r?`n$", ""
$g2 = $g2 -replace "`r?`n$", ""

$ws.Range("F2").Value = $f2 + "`n"
$ws.Range("G2").Value = $g2
